$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Image file names to add to column D (rows 2-11), one per kanji row.
$images = @(
    "n4_trip.jpeg",
    "n4_study.jpg",
    "n4_strong.png",
    "n4_teach.jpg",
    "n4_room.jpg",
    "n4_meeting.webp",
    "n4_company.webp",
    "n4_member.jpg",
    "n4_open.jpg",
    "n4_close.webp"
)

for ($i = 0; $i -lt $images.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $images[$i]
}

# Match the page setup (paper size / orientation) recorded for this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("F11").Select()
